# feat: add 2022-Q1 data
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Insert a new worksheet named "2022-Q1" right before the "总计" sheet.
#    This sheet will hold the per-fund holding detail for the new quarter,
#    following the same layout as the other quarterly sheets.
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$newSheet = $wb.Worksheets.Add($totalSheet)
$newSheet.Name = "2022-Q1"

# Header row
$newSheet.Range("B1").Value2 = "基金代码"
$newSheet.Range("C1").Value2 = "基金名称"
$newSheet.Range("D1").Value2 = "基金规模"
$newSheet.Range("E1").Value2 = "股票总仓位"
$newSheet.Range("F1").Value2 = "仓位占比"
$newSheet.Range("G1").Value2 = "持有市值(亿元)"
$newSheet.Range("H1").Value2 = "仓位排名"

$headerRange = $newSheet.Range("B1:H1")
$headerRange.Font.Bold = $true
$headerRange.Borders.LineStyle = 1
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160

# Fund holding detail rows
$data = @(
    @("010452", "广发瑞福精选混合A",         "16.29", "78.69", "3.37", "0.5490", 5),
    @("010453", "广发瑞福精选混合C",         "2.42",  "78.69", "3.37", "0.0816", 5),
    @("008443", "九泰动态策略灵活配置混合A", "0.24",  "64.82", "6.78", "0.0163", 2),
    @("007133", "嘉实长青竞争优势股票A",     "0.25",  "90.07", "4.55", "0.0114", 10),
    @("008444", "九泰动态策略灵活配置混合C", "0.13",  "64.82", "6.78", "0.0088", 2),
    @("005443", "国金量化多策略灵活配置混合", "0.51", "64.10", "0.91", "0.0046", 4),
    @("006346", "安信量化优选股票A",         "0.71",  "90.62", "0.61", "0.0043", 9),
    @("006347", "安信量化优选股票C",         "0.49",  "90.62", "0.61", "0.0030", 9),
    @("007134", "嘉实长青竞争优势股票C",     "0.02",  "90.07", "4.55", "0.0009", 10)
)

$row = 2
foreach ($item in $data) {
    $idxCell = $newSheet.Cells.Item($row, 1)
    $idxCell.Value2 = ($row - 2)
    $idxCell.Font.Bold = $true
    $idxCell.Borders.LineStyle = 1
    $idxCell.HorizontalAlignment = -4108
    $idxCell.VerticalAlignment = -4160

    $bCell = $newSheet.Cells.Item($row, 2)
    $bCell.NumberFormat = "@"
    $bCell.Value2 = $item[0]

    $newSheet.Cells.Item($row, 3).Value2 = $item[1]

    $dCell = $newSheet.Cells.Item($row, 4)
    $dCell.NumberFormat = "@"
    $dCell.Value2 = $item[2]

    $eCell = $newSheet.Cells.Item($row, 5)
    $eCell.NumberFormat = "@"
    $eCell.Value2 = $item[3]

    $fCell = $newSheet.Cells.Item($row, 6)
    $fCell.NumberFormat = "@"
    $fCell.Value2 = $item[4]

    $gCell = $newSheet.Cells.Item($row, 7)
    $gCell.NumberFormat = "@"
    $gCell.Value2 = $item[5]

    $newSheet.Cells.Item($row, 8).Value2 = $item[6]

    $row = $row + 1
}

# ---------------------------------------------------------------------------
# 2) Update the "总计" (summary) sheet: insert a new top data row for
#    2022-Q1 (9 holdings, 0.68 亿元) and push the existing rows down.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("总计")
$ws.Rows(2).Insert(-4121)
$ws.Range("A2:D2").ClearFormats()

$a2 = $ws.Range("A2")
$a2.Value2 = 0
$a2.Font.Bold = $true
$a2.Borders.LineStyle = 1
$a2.HorizontalAlignment = -4108
$a2.VerticalAlignment = -4160

$ws.Range("B2").Value2 = "2022-Q1"
$ws.Range("C2").Value2 = 9
$ws.Range("D2").Value2 = 0.68

# Renumber the index column (A) for the rows that were pushed down
for ($r = 3; $r -le 7; $r++) {
    $ws.Cells.Item($r, 1).Value2 = $r - 2
}
